# repull data, push all data, mean calculation
# Update column F (dSF) values for several rows to reflect repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F5").Value = -6
$ws.Range("F6").Value = -10
$ws.Range("F11").Value = -8
$ws.Range("F12").Value = -24
$ws.Range("F13").Value = -4
$ws.Range("F14").Value = -1
